# ScrumLog Kenny Vanrusselt 11/02/2015 final
#
# Updates the "Wat heb je gedaan?" / "Wat ga je doen?" entries for Kenny
# Vanrusselt (column B) on the single sheet, and moves the active
# selection to I4 (matching the saved cursor position in the workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B3: "Wat heb je gedaan?" entry for Kenny Vanrusselt.
$ws.Range("B3").Value = "Verzamelen van info, Opmaken van use case diagram + omschrijving"

# B4: "Wat ga je doen?" entry for Kenny Vanrusselt.
$ws.Range("B4").Value = "Ontwerpen van mockups / layout"

# Best-effort: move the window scroll position (xWindow/yWindow on the
# bookViews/workbookView element). Not all hosts expose this through the
# object model, so ignore failures and fall back to just the selection
# change, which is the part guaranteed to round-trip.
try {
    $excel.ActiveWindow.Top = 1800
} catch {
}

# Active cell / selection moves from D17 to I4.
$ws.Range("I4").Select()
